# Auto-generated PowerShell COM-interop script to apply the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows where only Price (D) and Volume(1h) (E) change ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.540.95'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.54%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.474.51'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.92%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.43%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '555.90'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.44%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '179.61'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.88%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.637'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.87%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.12%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.636'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.17%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.153'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.25%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '54.36'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.45%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000272'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.94%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.27'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.50%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.040.23'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.09%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '18.74'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.53%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.03'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.56%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '65.758.95'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.37%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.991'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.03%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '417.99'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.09%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.06'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.45%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '86.27'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.37%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.28'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.56%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.71'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +8.44%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.86'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -9.63%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.88'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.56%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.04'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.98%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.09'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.38%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '30.45'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.89%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.60'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.51%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '613.29'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -8.83%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.79'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.63%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.110'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.09%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '59.15'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.15%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0₃0797'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.10%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.340.93'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +10.00%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.383'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.90%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.31'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.27%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.45%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.84'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.75%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.56'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -9.17%  '

$ws.Range("E48").Value = '  -1.35%  '

$ws.Range("E49").Value = '  +2.40%  '

# --- Rows where Coin (B), Link (C), Price (D) and Volume(1h) (E) change (re-ranking) ---
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.485.80'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.09%  '

$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.121'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.20%  '

$ws.Range("B36").Value = 'Dai'
$ws.Range("C36").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.05%  '

$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.146'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +9.63%  '

$ws.Range("B38").Value = 'InjectiveProtocol'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '37.65'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.71%  '

$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.28'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.30%  '

$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0416'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.58%  '

$ws.Range("B50").Value = 'THORChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.47'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.09%  '

$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '137.89'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.26%  '
